$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 61 data (K Closest Points to Origin / 973) ---
$ws.Range("A61").Value = "K Closest Points to Origin"
$ws.Range("B61").Value = "Array"
$ws.Range("C61").Value = "No"
$ws.Range("D61").Value = "Yes"
$ws.Range("E61").Value = "Medium"
$ws.Range("F61").Value = "Easy"

$g61 = $ws.Range("G61")
$link = $ws.Hyperlinks.Add($g61, "973%20-%20K%20Closest%20Points%20to%20Origin")
$link.TextToDisplay = "973 - K Closest Points to Origin"
$g61.Style = "Hyperlink"

# --- Extend data validations to include row 61 (recreate in original order) ---
$ws.Range("E2:F60").Validation.Delete()
$ws.Range("E2:F61").Validation.Add(3, 1, 1, '"Easy, Medium, Hard"')

$cVal = $ws.Range("C2:C60")
$cVal.Validation.Delete()
$cValNew = $ws.Range("C2:C61")
$cValNew.Validation.Add(3, 1, 1, '"Yes, No"')
$cValNew.Validation.IgnoreBlank = $false

$ws.Range("D2:D60").Validation.Delete()
$ws.Range("D2:D61").Validation.Add(3, 1, 1, '"Yes, No"')

$ws.Range("B2:B60").Validation.Delete()
$ws.Range("B2:B61").Validation.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap, Class Design"')

# --- Extend conditional formatting to include new row. The original file
# has ONE rule-group (3 cfRules) shared across the two areas "D2:G8" and
# "D9:F60" (rendered as sqref="D2:G8 D9:F60"). The COM model here only
# supports adding a rule to a single contiguous area at a time, and
# deleting from either area's FormatConditions collection deletes from the
# whole shared group. So: delete the shared group once, then recreate the
# same 3 rules (same formulas/colors/priorities) separately on "D2:G8" and
# on the row-extended "D9:F61" area. ---
$cf = $ws.Range("D2:G8,D9:F60").FormatConditions
while ($cf.Count -gt 0) {
  $cf.Item(1).Delete()
}

$areaTop = $ws.Range("D2:G8")

$fcHardTop = $areaTop.FormatConditions.Add(1, 3, '="Hard"')
$fcHardTop.Interior.Color = 255
$fcHardTop.Priority = 7

$fcMediumTop = $areaTop.FormatConditions.Add(1, 3, '="Medium"')
$fcMediumTop.Interior.Color = 5296274
$fcMediumTop.Priority = 8

$fcEasyTop = $areaTop.FormatConditions.Add(1, 3, '="Easy"')
$fcEasyTop.Interior.Color = 65280
$fcEasyTop.Priority = 9

$areaBottom = $ws.Range("D9:F61")

$fcHardBottom = $areaBottom.FormatConditions.Add(1, 3, '="Hard"')
$fcHardBottom.Interior.Color = 255
$fcHardBottom.Priority = 7

$fcMediumBottom = $areaBottom.FormatConditions.Add(1, 3, '="Medium"')
$fcMediumBottom.Interior.Color = 5296274
$fcMediumBottom.Priority = 8

$fcEasyBottom = $areaBottom.FormatConditions.Add(1, 3, '="Easy"')
$fcEasyBottom.Interior.Color = 65280
$fcEasyBottom.Priority = 9

# --- Update used-range dimension & current selection to match the edit ---
$ws.Range("N46").Select()

Write-Output "done"
